# Insert a new column before F: shifts old District values from F to G,
# creating an empty F column for the new "Address" data (matches the diff,
# which shows F2 gaining "Address" while the former F column content
# (District) reappears unchanged under the new G column).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("F:F").Insert()

# Header
$ws.Range("F2").Value = "Address"

# Address values (school name/location extracted from column B), row by row.
# Rows 17 and 34 are intentionally left blank -- in the source data those rows
# hold only a fragment (continuation) of the text from the row above, so the
# diff leaves their new Address cell empty too.
$ws.Range("F3").Value = "G H S Bijjahalli"
$ws.Range("F4").Value = "G H S MadbalMagadi"
$ws.Range("F5").Value = "G H S Srigiripura Magadi"
$ws.Range("F6").Value = "Govt. High School Arasana KunteMagadi"
$ws.Range("F7").Value = "G H S Doddamara LavadrKanakapura"
$ws.Range("F8").Value = "G H S Halasurukanakpura"
$ws.Range("F9").Value = "G H S ThigalarahosahalliKanakapura"
$ws.Range("F10").Value = "Govt. Adarsha Vidyalaya KodihalliKanakapura"
$ws.Range("F11").Value = "G H S BannimukkodluKanakapura"
$ws.Range("F12").Value = "Govt. High School AlanathaKanakapura"
$ws.Range("F13").Value = "G H S MangalgiHumnabad"
$ws.Range("F14").Value = "Govt. High School HarkoodBasavakalyan"
$ws.Range("F15").Value = "Govt. High School MadakattiBhalki"
$ws.Range("F16").Value = "G H S DanagapurBhalki"
# F17 intentionally left blank
$ws.Range("F18").Value = "Govt. Girls High SchoolThana KushnvorAurad"
$ws.Range("F19").Value = "Govt. High SchoolJanatangar HudgiHumnabad"
$ws.Range("F20").Value = "G H S LadhaBhalki"
$ws.Range("F21").Value = "Govt. Hr. P. S Nealg. Bhalki"
$ws.Range("F22").Value = "Govt. High School GourBasavakalyan"
$ws.Range("F23").Value = "G H S Hallikhed(K) Humnabad"
$ws.Range("F24").Value = "Govt. H P S MorambiBhalki"
$ws.Range("F25").Value = "Govt. High School Manahalli"
$ws.Range("F26").Value = "Govt. High SchoolBadalgaonAurad"
$ws.Range("F27").Value = "G H S Ranjol Kheni"
$ws.Range("F28").Value = "G Hr P S BolegaonBhalki"
$ws.Range("F29").Value = "Joyti High School KushnoorAurad"
$ws.Range("F30").Value = "Govt. High SchoolNagmarpalliAurad"
$ws.Range("F31").Value = "G Hr P S MadakattiBhalki"
$ws.Range("F32").Value = "Sarvodaya Girl high schoolHumnabad"
$ws.Range("F33").Value = "Govt. High School Gadlegaon(B)"
# F34 intentionally left blank
$ws.Range("F35").Value = "Govt. High School Dhannur (K)Basavakalyan"
$ws.Range("F36").Value = "G H S Gadagi"
$ws.Range("F37").Value = "Govt. High School LakhangaonBhalki"
$ws.Range("F38").Value = "G H P S KoudgaonAurad"
